# Update the "Fitness" values in column C (Sheet1) per the target dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C3").Value = 12126
$ws.Range("C4:C5").Value = 9929
$ws.Range("C6:C11").Value = 9531
$ws.Range("C12").Value = 9468
$ws.Range("C13").Value = 9169
$ws.Range("C14:C15").Value = 9074
$ws.Range("C16:C21").Value = 9027
$ws.Range("C22:C26").Value = 8998
$ws.Range("C27:C31").Value = 8786
$ws.Range("C32:C42").Value = 8511
$ws.Range("C43:C51").Value = 7939
$ws.Range("C52:C56").Value = 7925
$ws.Range("C57:C65").Value = 7678
$ws.Range("C66:C97").Value = 7594
$ws.Range("C190:C252").Value = 7569
